$wb = $excel.ActiveWorkbook

# --- Sheet: Significant Components ---
$ws = $wb.Worksheets.Item('Significant Components')
$ws.Range('C2').Value = '[''QESL'' ''PPUNIT'' ''QMOHO'' ''QEDLESHI'' ''QHISPC'' ''QNOHLTH'' ''QEXTRCT'']'
$ws.Range('C3').Value = '[''PERCAP'' ''QRICH'' ''MDHSEVAL'']'
$ws.Range('C5').Value = '[''QRENTER'' ''QNOAUTO'' ''QPOVTY'']'
$ws.Range('C6').Value = '[''QAGEDEP'' ''QFEMALE'' ''QFEMLBR'']'

# --- Sheet: Loading Factors ---
$ws = $wb.Worksheets.Item('Loading Factors')
$ws.Range('A2').Value = 'QESL'
$ws.Range('B2').Value = 0.7463154738714127
$ws.Range('C2').Value = 0.1383746688021929
$ws.Range('D2').Value = -0.06828868531064654
$ws.Range('E2').Value = 0.2025402315337327
$ws.Range('F2').Value = -0.04405009710752585
$ws.Range('G2').Value = 0.09086866093976122
$ws.Range('A3').Value = 'PPUNIT'
$ws.Range('B3').Value = 0.5835173019418969
$ws.Range('C3').Value = -0.03689747254769649
$ws.Range('D3').Value = -0.1272499913564337
$ws.Range('E3').Value = -0.4204908240952777
$ws.Range('F3').Value = 0.1006093523631547
$ws.Range('G3').Value = 0.130977512453249
$ws.Range('A4').Value = 'QMOHO'
$ws.Range('B4').Value = 0.4515387097126257
$ws.Range('C4').Value = 0.1342118373807977
$ws.Range('D4').Value = 0.03533061792857509
$ws.Range('E4').Value = -0.1630772623581569
$ws.Range('F4').Value = -0.01675659521983348
$ws.Range('G4').Value = 0.001862665364592792
$ws.Range('A5').Value = 'QEDLESHI'
$ws.Range('B5').Value = 0.8352779339073224
$ws.Range('C5').Value = 0.1442623953550782
$ws.Range('D5').Value = 0.002317829691818668
$ws.Range('E5').Value = 0.1530479582561485
$ws.Range('F5').Value = -0.02604027049597342
$ws.Range('G5').Value = 0.1765466381914305
$ws.Range('A6').Value = 'QHISPC'
$ws.Range('B6').Value = 0.7423692396666269
$ws.Range('C6').Value = 0.3098526168554369
$ws.Range('D6').Value = -0.1486952812711276
$ws.Range('E6').Value = 0.09400583957306605
$ws.Range('F6').Value = -0.04885413376749048
$ws.Range('G6').Value = 0.2354826922092088
$ws.Range('A7').Value = 'QNOHLTH'
$ws.Range('B7').Value = 0.6038940751127572
$ws.Range('C7').Value = 0.3577588845444675
$ws.Range('D7').Value = -0.07977624897329463
$ws.Range('E7').Value = 0.2571839833933484
$ws.Range('F7').Value = -0.09534315052672397
$ws.Range('G7').Value = 0.1820152857384868
$ws.Range('A8').Value = 'QEXTRCT'
$ws.Range('B8').Value = 0.7662017861373808
$ws.Range('C8').Value = 0.1246192445311004
$ws.Range('D8').Value = 0.05703178003480995
$ws.Range('E8').Value = 0.1061640337901088
$ws.Range('F8').Value = -0.1262460497666379
$ws.Range('G8').Value = 0.00427583026953045
$ws.Range('B9').Value = 0.4141283943230038
$ws.Range('C9').Value = 0.6768765468547735
$ws.Range('D9').Value = -0.208960726327434
$ws.Range('E9').Value = 0.189163545024027
$ws.Range('F9').Value = 0.09876521159871486
$ws.Range('G9').Value = 0.1797756325458044
$ws.Range('A10').Value = 'QRICH'
$ws.Range('B10').Value = 0.1786760063468116
$ws.Range('C10').Value = 0.8099499573541392
$ws.Range('D10').Value = -0.1171137254553127
$ws.Range('E10').Value = 0.3686627495693146
$ws.Range('F10').Value = -0.02002121773364113
$ws.Range('G10').Value = 0.109511635490191
$ws.Range('A11').Value = 'MDHSEVAL'
$ws.Range('B11').Value = 0.353757620337989
$ws.Range('C11').Value = 0.728708382215542
$ws.Range('D11').Value = -0.08538345655936434
$ws.Range('E11').Value = -0.04214449505572827
$ws.Range('F11').Value = 0.04384743220829919
$ws.Range('G11').Value = 0.1201375945681547
$ws.Range('A12').Value = 'QRENTER'
$ws.Range('B12').Value = -0.06964050153989358
$ws.Range('C12').Value = 0.2243052317807107
$ws.Range('D12').Value = -0.4250942622562807
$ws.Range('E12').Value = 0.7486728906150069
$ws.Range('F12').Value = -0.08020597109228203
$ws.Range('G12').Value = 0.06481865033512572
$ws.Range('A13').Value = 'MEDAGE'
$ws.Range('B13').Value = -0.1933922661302453
$ws.Range('C13').Value = -0.2556005570341998
$ws.Range('D13').Value = 0.7735816192592738
$ws.Range('E13').Value = -0.2675924831286911
$ws.Range('F13').Value = -0.05058706228423991
$ws.Range('G13').Value = -0.1226097589316303
$ws.Range('B14').Value = 0.03196500328660701
$ws.Range('C14').Value = -0.06645709694710349
$ws.Range('D14').Value = 0.6956263850825187
$ws.Range('E14').Value = -0.04457261547133459
$ws.Range('F14').Value = 0.6123584634842314
$ws.Range('G14').Value = -0.04113476791561829
$ws.Range('B15').Value = 0.02114278701105845
$ws.Range('C15').Value = -0.03947992616860568
$ws.Range('D15').Value = 0.7957449187436122
$ws.Range('E15').Value = -0.08782706873952935
$ws.Range('F15').Value = 0.06483515366448697
$ws.Range('G15').Value = 0.04245536209771213
$ws.Range('A16').Value = 'QNOAUTO'
$ws.Range('B16').Value = 0.05709370144763014
$ws.Range('C16').Value = 0.04107575786895225
$ws.Range('D16').Value = 0.01144308742084319
$ws.Range('E16').Value = 0.6493597865939288
$ws.Range('F16').Value = 0.0149432863733358
$ws.Range('G16').Value = 0.05784573758531825
$ws.Range('A17').Value = 'QPOVTY'
$ws.Range('B17').Value = 0.1957216580403172
$ws.Range('C17').Value = 0.1330977692004899
$ws.Range('D17').Value = -0.2320824202469326
$ws.Range('E17').Value = 0.5572080744446453
$ws.Range('F17').Value = 0.08212745736725391
$ws.Range('G17').Value = 0.1166838084313597
$ws.Range('A18').Value = 'QFEMALE'
$ws.Range('B18').Value = 0.02896968068441856
$ws.Range('C18').Value = 0.001264420953605554
$ws.Range('D18').Value = 0.1484947721477394
$ws.Range('E18').Value = 0.04764862560814205
$ws.Range('F18').Value = 0.9303837596979473
$ws.Range('G18').Value = -0.0001998259267016299
$ws.Range('A19').Value = 'QFEMLBR'
$ws.Range('B19').Value = -0.2164852227571602
$ws.Range('C19').Value = 0.07008131099298208
$ws.Range('D19').Value = -0.05054873711277632
$ws.Range('E19').Value = -0.01544327837111355
$ws.Range('F19').Value = 0.6764703176186838
$ws.Range('G19').Value = 0.1327983158183232
$ws.Range('B20').Value = 0.1389414203490118
$ws.Range('C20').Value = 0.2227513400325371
$ws.Range('D20').Value = -0.1031614102698809
$ws.Range('E20').Value = 0.3534990123545979
$ws.Range('F20').Value = -0.02460918471860097
$ws.Range('G20').Value = 0.4969878288289207
$ws.Range('B21').Value = 0.3059368164415499
$ws.Range('C21').Value = 0.1540730302632215
$ws.Range('D21').Value = 0.02452183562253974
$ws.Range('E21').Value = 0.02957443205099398
$ws.Range('F21').Value = 0.1777865111549735
$ws.Range('G21').Value = 0.9169224125282442

# --- Sheet: All Refactor Variances ---
$ws = $wb.Worksheets.Item('All Refactor Variances')
$ws.Range('B2').Value = 4.095984283209816
$ws.Range('C2').Value = 2.853394873491776
$ws.Range('D2').Value = 2.240850447181389
$ws.Range('E2').Value = 1.887013627171163
$ws.Range('F2').Value = 1.826686531450328
$ws.Range('G2').Value = 1.571823091120927
$ws.Range('H2').Value = 0.9800951980275844
$ws.Range('I2').Value = 3.876385829778928
$ws.Range('J2').Value = 2.1633517207822
$ws.Range('K2').Value = 2.104589594732921
$ws.Range('L2').Value = 2.029886881000504
$ws.Range('M2').Value = 1.803101094454233
$ws.Range('N2').Value = 1.348998893620019
$ws.Range('B3').Value = 0.1517031216003636
$ws.Range('C3').Value = 0.1056812916108065
$ws.Range('D3').Value = 0.0829944610067181
$ws.Range('E3').Value = 0.06988939359893195
$ws.Range('F3').Value = 0.06765505672038251
$ws.Range('G3').Value = 0.05821567004151582
$ws.Range('H3').Value = 0.03629982214916979
$ws.Range('I3').Value = 0.1938192914889464
$ws.Range('J3').Value = 0.10816758603911
$ws.Range('K3').Value = 0.105229479736646
$ws.Range('L3').Value = 0.1014943440500252
$ws.Range('M3').Value = 0.09015505472271164
$ws.Range('N3').Value = 0.06744994468100093
$ws.Range('B4').Value = 0.1517031216003636
$ws.Range('C4').Value = 0.2573844132111701
$ws.Range('D4').Value = 0.3403788742178881
$ws.Range('E4').Value = 0.4102682678168201
$ws.Range('F4').Value = 0.4779233245372027
$ws.Range('G4').Value = 0.5361389945787185
$ws.Range('H4').Value = 0.5724388167278883
$ws.Range('I4').Value = 0.1938192914889464
$ws.Range('J4').Value = 0.3019868775280564
$ws.Range('K4').Value = 0.4072163572647025
$ws.Range('L4').Value = 0.5087107013147276
$ws.Range('M4').Value = 0.5988657560374393
$ws.Range('N4').Value = 0.6663157007184403
$ws.Range('B5').Value = 0.2650119404332366
$ws.Range('C5').Value = 0.1846158725134858
$ws.Range('D5').Value = 0.1449839853298592
$ws.Range('E5').Value = 0.1220905912677725
$ws.Range('F5').Value = 0.1181874022923618
$ws.Range('G5').Value = 0.1016976283584014
$ws.Range('H5').Value = 0.06341257980488262
$ws.Range('I5').Value = 0.2908820717866396
$ws.Range('J5').Value = 0.1623368411137253
$ws.Range('K5').Value = 0.1579273602936035
$ws.Range('L5').Value = 0.1523217056728382
$ws.Range('M5').Value = 0.1353038126304152
$ws.Range('N5').Value = 0.1012282085027781

# --- Sheet: Final Variances ---
$ws = $wb.Worksheets.Item('Final Variances')
$ws.Range('B2').Value = 3.876385829778928
$ws.Range('C2').Value = 2.1633517207822
$ws.Range('D2').Value = 2.104589594732921
$ws.Range('E2').Value = 2.029886881000504
$ws.Range('F2').Value = 1.803101094454233
$ws.Range('G2').Value = 1.348998893620019
$ws.Range('B3').Value = 0.1938192914889464
$ws.Range('C3').Value = 0.10816758603911
$ws.Range('D3').Value = 0.105229479736646
$ws.Range('E3').Value = 0.1014943440500252
$ws.Range('F3').Value = 0.09015505472271164
$ws.Range('G3').Value = 0.06744994468100093
$ws.Range('B4').Value = 0.1938192914889464
$ws.Range('C4').Value = 0.3019868775280564
$ws.Range('D4').Value = 0.4072163572647025
$ws.Range('E4').Value = 0.5087107013147276
$ws.Range('F4').Value = 0.5988657560374393
$ws.Range('G4').Value = 0.6663157007184403
$ws.Range('B5').Value = 0.2908820717866396
$ws.Range('C5').Value = 0.1623368411137253
$ws.Range('D5').Value = 0.1579273602936035
$ws.Range('E5').Value = 0.1523217056728382
$ws.Range('F5').Value = 0.1353038126304152
$ws.Range('G5').Value = 0.1012282085027781

# --- Sheet: Included and Excluded ---
$ws = $wb.Worksheets.Item('Included and Excluded')
$ws.Range('B2').Value = '[[''QESL'', ''PPUNIT'', ''QMOHO'', ''QEDLESHI'', ''QHISPC'', ''QNOHLTH'', ''QEXTRCT'', ''PERCAP'', ''QRICH'', ''MDHSEVAL'', ''MEDAGE'', ''QAGEDEP'', ''QSSBEN'', ''QRENTER'', ''QNOAUTO'', ''QPOVTY'', ''QFEMALE'', ''QFEMLBR'', ''QFAM'', ''QFHH'']]'

